$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Minor timestamp precision correction on the previous row (row 7), as
# produced by the scheduled task's re-computation.
$ws.Cells.Item(7, 1).Value = 45864.45857037037

# Append a new row of sensor data (row 8), matching the style of row 7.
$row = 8

$ws.Cells.Item($row, 1).Value = 45864.5003065265
$ws.Cells.Item($row, 1).NumberFormat = $ws.Cells.Item($row - 1, 1).NumberFormat

$ws.Cells.Item($row, 2).Value = 2025
$ws.Cells.Item($row, 3).Value = 30
$ws.Cells.Item($row, 4).Value = 19.94
$ws.Cells.Item($row, 5).Value = 71.17
$ws.Cells.Item($row, 6).Value = 625.78
$ws.Cells.Item($row, 7).Value = 11.77
$ws.Cells.Item($row, 8).Value = "SE"
$ws.Cells.Item($row, 9).Value = 0
$ws.Cells.Item($row, 10).Value = "12:00:26"
